$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns K:P on row 1
$ws.Range("K1").Value = "From1"
$ws.Range("L1").Value = "To1"
$ws.Range("M1").Value = "From2"
$ws.Range("N1").Value = "To2"
$ws.Range("O1").Value = "From3"
$ws.Range("P1").Value = "To3"

# New data cells for columns K:P on row 4 (MultiCity row)
$ws.Range("K4").Value = "Mumbai"
$ws.Range("L4").Value = "Pune"
$ws.Range("M4").Value = "Pune"
$ws.Range("N4").Value = "Mumbai"
$ws.Range("O4").Value = "Mumbai"
$ws.Range("P4").Value = "Delhi"

# Update the selected cell to match the target view state
$ws.Range("N9").Select()
